$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 12.0491
$ws.Range("E7").Value = 11.9606
$ws.Range("B8").Value = 4.947300000000003
$ws.Range("A12").Value = -22.75470000000001
$ws.Range("B12").Value = 5.868300000000001
$ws.Range("B14").Value = 8.829700000000004
$ws.Range("E19").Value = 12.98519999999999
$ws.Range("E21").Value = 12.7924
$ws.Range("B22").Value = 4.638000000000008
$ws.Range("E24").Value = 12.41159999999999
